$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''41.532.82'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.59%  '
$ws.Range('D3').Value = '''2.482.47'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.36%  '
$ws.Range('D4').Value = '''0.997'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.27%  '
$ws.Range('D5').Value = '''312.11'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.80%  '
$ws.Range('D6').Value = '''93.08'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.01%  '
$ws.Range('D7').Value = '''0.545'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.54%  '
$ws.Range('D8').Value = '''0.999'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.34%  '
$ws.Range('D9').Value = '''0.497'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.90%  '
$ws.Range('D10').Value = '''32.72'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.04%  '
$ws.Range('D11').Value = '''0.0783'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.49%  '
$ws.Range('E12').Value = '  +2.26%  '
$ws.Range('D13').Value = '''2.867.32'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.28%  '
$ws.Range('D14').Value = '''6.87'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.95%  '
$ws.Range('D15').Value = '''15.47'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +7.35%  '
$ws.Range('D16').Value = '''2.480.18'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.49%  '
$ws.Range('D17').Value = '''0.755'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.93%  '
$ws.Range('D18').Value = '''41.676.13'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.81%  '
$ws.Range('E19').Value = '  -0.23%  '
$ws.Range('D20').Value = '''0.0₃0921'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.14%  '
$ws.Range('D21').Value = '''70.75'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +5.85%  '
$ws.Range('D22').Value = '''11.17'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.07%  '
$ws.Range('D23').Value = '''235.60'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.43%  '
$ws.Range('D24').Value = '''2.70'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.06%  '
$ws.Range('E25').Value = '  -0.13%  '
$ws.Range('E26').Value = '  -0.76%  '
$ws.Range('E27').Value = '  +1.18%  '
$ws.Range('D28').Value = '''2.24'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.49%  '
$ws.Range('D29').Value = '''9.63'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.27%  '
$ws.Range('D30').Value = '''36.24'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.39%  '
$ws.Range('D31').Value = '''153.84'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.52%  '
$ws.Range('E32').Value = '  -3.08%  '
$ws.Range('B33').Value = 'WEMIXToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D33').Value = '''2.58'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.86%  '
$ws.Range('B34').Value = 'Celestia'
$ws.Range('C34').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D34').Value = '''18.16'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +5.81%  '
$ws.Range('E35').Value = '  +0.82%  '
$ws.Range('D36').Value = '''2.47'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.63%  '
$ws.Range('D37').Value = '''2.95'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.17%  '
$ws.Range('E38').Value = '  -3.55%  '
$ws.Range('E39').Value = '  +0.44%  '
$ws.Range('E40').Value = '  -0.16%  '
$ws.Range('E41').Value = '  -0.36%  '
$ws.Range('D42').Value = '''1.00'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.03%  '
$ws.Range('D43').Value = '''19.80'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -6.22%  '
$ws.Range('D44').Value = '''1.952.10'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.75%  '
$ws.Range('D45').Value = '''0.0284'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.06%  '
$ws.Range('D46').Value = '''2.96'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.65%  '
$ws.Range('E47').Value = '  +1.45%  '
$ws.Range('D48').Value = '''2.730.95'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.48%  '
$ws.Range('D49').Value = '''96.20'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.05%  '
$ws.Range('D50').Value = '''67.07'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.69%  '
$ws.Range('D51').Value = '''73.26'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.98%  '
